# Updates the "Estado de Cuenta" detail rows (Hoja1!C16:G51) to the new
# database snapshot: replaces the previous SEBASTIAN JOSE MARTELO ESQUIVEL /
# SANDRA MELINA GARNICA FARAK / YEAN CARLOS MEZA PITALUA / LUIS EDUARDO
# BOLIVAR MENDOZA block with the refreshed periods/values, inserting the new
# "parte 1" period (1902) for LUIS EDUARDO BOLIVAR MENDOZA at row 48 and
# shifting the remaining rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, C=N Doc Trabajador, D=Nombre Trabajador, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    @(16, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1607", 25774, 1000000),
    @(17, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1608", 25774, 1000000),
    @(18, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1609", 25774, 1000000),
    @(19, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1610", 25774, 1000000),
    @(20, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1611", 25774, 1000000),
    @(21, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1612", 25774, 1000000),
    @(22, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1701", 25774, 1000000),
    @(23, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1702", 25774, 1000000),
    @(24, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1703", 25774, 1000000),
    @(25, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1704", 25774, 1000000),
    @(26, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1705", 25774, 1000000),
    @(27, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1706", 25774, 1000000),
    @(28, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1707", 25774, 1000000),
    @(29, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1708", 25774, 1000000),
    @(30, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1709", 25774, 1000000),
    @(31, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1710", 25774, 1000000),
    @(32, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1711", 25774, 1000000),
    @(33, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1712", 25774, 1000000),
    @(34, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1801", 25774, 1000000),
    @(35, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1802", 25774, 1000000),
    @(36, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1803", 25774, 1000000),
    @(37, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1804", 25774, 1000000),
    @(38, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1805", 25774, 1000000),
    @(39, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1806", 25774, 1000000),
    @(40, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1807", 25774, 1000000),
    @(41, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1808", 25774, 1000000),
    @(42, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1809", 31249, 1000000),
    @(43, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1810", 31249, 1000000),
    @(44, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1811", 31249, 1000000),
    @(45, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1812", 31249, 1000000),
    @(46, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1901", 31249, 1000000),
    @(47, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1902", 31249, 1000000),
    @(48, "1047474222", "LUIS EDUARDO BOLIVAR MENDOZA", "1902", 4417, 828116),
    @(49, "1072260135", "SEBASTIAN JOSE MARTELO ESQUIVEL", "1903", 31249, 1000000),
    @(50, "43611573", "SANDRA MELINA GARNICA FARAK", "1911", 6625, 877803),
    @(51, "1143365578", "YEAN CARLOS MEZA PITALUA", "1911", 6625, 877803)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
